# Applies "atualizei dados da bibi e add" update to Resumo_por_Cliente sheet
# (previsao_retorno.xlsx) - refreshed "meses sem comprar" counters plus a
# handful of recomputed probability/date/count fields for a few clients.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

# --- Row 14 (id_cliente 2095): total_compras_historico + date window shifted ---
$ws.Range("E14").Value = 28
$ws.Range("H14").Value = 45852.53395833333
$ws.Range("I14").Value = 45914.53395833333

# --- Row 28 (id_cliente 3355): probabilities, counts, padrao_compra & dates recomputed ---
$ws.Range("B28").Value = 0.42
$ws.Range("D28").Value = 0.5
$ws.Range("E28").Value = 11
$ws.Range("F28").Value = 0.5
$ws.Range("G28").Value = "1x a cada 4 meses - irregular"
$ws.Range("H28").Value = 45852.67643518518
$ws.Range("I28").Value = 45975.67643518518

# --- Row 115 (id_cliente 28458): total_compras_historico + date window shifted ---
$ws.Range("E115").Value = 16612
$ws.Range("H115").Value = 45852.67851851852
$ws.Range("I115").Value = 45853.67851851852

# --- "situacao" (column J) refreshed "meses sem comprar" counters ---
$ws.Range("J5").Value   = "INATIVO - 15.4 meses sem comprar"
$ws.Range("J6").Value   = "INATIVO - 16.9 meses sem comprar"
$ws.Range("J16").Value  = "INATIVO - 40.4 meses sem comprar"
$ws.Range("J24").Value  = "INATIVO - 38.4 meses sem comprar"
$ws.Range("J25").Value  = "INATIVO - 0.4 meses sem comprar"
$ws.Range("J31").Value  = "INATIVO - 7.4 meses sem comprar"
$ws.Range("J39").Value  = "INATIVO - 32.8 meses sem comprar"
$ws.Range("J46").Value  = "INATIVO - 6.7 meses sem comprar"
$ws.Range("J47").Value  = "INATIVO - 16.5 meses sem comprar"
$ws.Range("J51").Value  = "INATIVO - 8.1 meses sem comprar"
$ws.Range("J66").Value  = "INATIVO - 28.4 meses sem comprar"
$ws.Range("J70").Value  = "INATIVO - 11.9 meses sem comprar"
$ws.Range("J77").Value  = "INATIVO - 8.3 meses sem comprar"
$ws.Range("J80").Value  = "INATIVO - 6.7 meses sem comprar"
$ws.Range("J89").Value  = "INATIVO - 15.6 meses sem comprar"
$ws.Range("J91").Value  = "INATIVO - 12.3 meses sem comprar"
$ws.Range("J92").Value  = "INATIVO - 11.7 meses sem comprar"
$ws.Range("J99").Value  = "INATIVO - 33.5 meses sem comprar"
$ws.Range("J103").Value = "INATIVO - 37.8 meses sem comprar"
$ws.Range("J105").Value = "INATIVO - 15.1 meses sem comprar"
